{"js": "// Replace the date line and each \"A\u00d7B=\" expression in the practice table\n// with the new values, as described by the diff. Every source string in\n// this document is unique, so a simple search-and-replace keyed on the\n// exact old text is unambiguous and preserves all existing run formatting\n// (font, size, etc.) because we replace only the matched range's text.\n\nconst replacements = [\n  [\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"],\n  [\"189\u00d79=\", \"152\u00d74=\"],\n  [\"267\u00d74=\", \"288\u00d79=\"],\n  [\"761\u00d76=\", \"220\u00d72=\"],\n  [\"973\u00d76=\", \"249\u00d73=\"],\n  [\"283\u00d75=\", \"387\u00d72=\"],\n  [\"991\u00d79=\", \"470\u00d76=\"],\n  [\"958\u00d76=\", \"851\u00d77=\"],\n  [\"698\u00d74=\", \"702\u00d75=\"],\n  [\"901\u00d72=\", \"654\u00d75=\"],\n  [\"713\u00d72=\", \"647\u00d75=\"],\n  [\"911\u00d78=\", \"379\u00d73=\"],\n  [\"729\u00d76=\", \"106\u00d77=\"],\n  [\"131\u00d72=\", \"641\u00d73=\"],\n  [\"231\u00d76=\", \"120\u00d76=\"],\n  [\"769\u00d73=\", \"370\u00d77=\"],\n  [\"463\u00d74=\", \"448\u00d78=\"],\n  [\"598\u00d75=\", \"541\u00d77=\"],\n  [\"969\u00d78=\", \"532\u00d77=\"],\n  [\"444\u00d79=\", \"393\u00d73=\"],\n  [\"910\u00d75=\", \"947\u00d78=\"],\n  [\"492\u00d77=\", \"406\u00d73=\"],\n  [\"429\u00d79=\", \"157\u00d72=\"],\n  [\"396\u00d73=\", \"704\u00d76=\"],\n  [\"555\u00d74=\", \"581\u00d79=\"],\n  [\"295\u00d75=\", \"926\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00d7B=\" expression in the practice table\n# with the new values, as described by the diff. Every source string in\n# this document is unique, so a simple Find/Replace keyed on the exact old\n# text is unambiguous and preserves existing run formatting (font, size,\n# etc.) because Word's Find & Replace substitutes only the matched text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-09-09 Tuesday\", \"2025-09-10 Wednesday\"),\n  @(\"189\u00d79=\", \"152\u00d74=\"),\n  @(\"267\u00d74=\", \"288\u00d79=\"),\n  @(\"761\u00d76=\", \"220\u00d72=\"),\n  @(\"973\u00d76=\", \"249\u00d73=\"),\n  @(\"283\u00d75=\", \"387\u00d72=\"),\n  @(\"991\u00d79=\", \"470\u00d76=\"),\n  @(\"958\u00d76=\", \"851\u00d77=\"),\n  @(\"698\u00d74=\", \"702\u00d75=\"),\n  @(\"901\u00d72=\", \"654\u00d75=\"),\n  @(\"713\u00d72=\", \"647\u00d75=\"),\n  @(\"911\u00d78=\", \"379\u00d73=\"),\n  @(\"729\u00d76=\", \"106\u00d77=\"),\n  @(\"131\u00d72=\", \"641\u00d73=\"),\n  @(\"231\u00d76=\", \"120\u00d76=\"),\n  @(\"769\u00d73=\", \"370\u00d77=\"),\n  @(\"463\u00d74=\", \"448\u00d78=\"),\n  @(\"598\u00d75=\", \"541\u00d77=\"),\n  @(\"969\u00d78=\", \"532\u00d77=\"),\n  @(\"444\u00d79=\", \"393\u00d73=\"),\n  @(\"910\u00d75=\", \"947\u00d78=\"),\n  @(\"492\u00d77=\", \"406\u00d73=\"),\n  @(\"429\u00d79=\", \"157\u00d72=\"),\n  @(\"396\u00d73=\", \"704\u00d76=\"),\n  @(\"555\u00d74=\", \"581\u00d79=\"),\n  @(\"295\u00d75=\", \"926\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
